# Sync file from Google Drive
# Refresh the live NextBus arrival snapshots (ETA timestamps, minutes-to-arrival,
# and a handful of bus-type / monitored flags) on each of the three sheets to
# match the latest pull from the source feed.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("NextBus1")
$ws.Range("F2").Value = 45688.43534722222
$ws.Range("O2").Value = 15
$ws.Range("F3").Value = 45688.43188657407
$ws.Range("L3").Value = "DD"
$ws.Range("O3").Value = 10
$ws.Range("F4").Value = 45688.42649305556
$ws.Range("L4").Value = "DD"
$ws.Range("O4").Value = 2
$ws.Range("F5").Value = 45688.43152777778
$ws.Range("O5").Value = 10
$ws.Range("F6").Value = 45688.42457175926
$ws.Range("O6").Value = 0
$ws.Range("F7").Value = 45688.42671296297
$ws.Range("F8").Value = 45688.43047453704
$ws.Range("L8").Value = "DD"
$ws.Range("O8").Value = 8
$ws.Range("F9").Value = 45688.42724537037
$ws.Range("O9").Value = 3
$ws.Range("F10").Value = 45688.43097222222
$ws.Range("O10").Value = 9
$ws.Range("F11").Value = 45688.42685185185
$ws.Range("O11").Value = 3
$ws.Range("F12").Value = 45688.4240625
$ws.Range("O12").Value = 0
$ws.Range("F13").Value = 45688.43243055556
$ws.Range("O13").Value = 11
$ws.Range("F14").Value = 45688.42769675926
$ws.Range("O14").Value = 4
$ws.Range("F15").Value = 45688.42576388889
$ws.Range("O15").Value = 1

$ws = $wb.Worksheets.Item("NextBus2")
$ws.Range("F2").Value = 45688.44231481481
$ws.Range("O2").Value = 25
$ws.Range("F3").Value = 45688.43688657408
$ws.Range("O3").Value = 17
$ws.Range("F4").Value = 45688.4350462963
$ws.Range("L4").Value = "SD"
$ws.Range("O4").Value = 15
$ws.Range("F5").Value = 45688.44203703704
$ws.Range("O5").Value = 25
$ws.Range("F6").Value = 45688.43158564815
$ws.Range("O6").Value = 10
$ws.Range("F7").Value = 45688.43606481481
$ws.Range("J7").Value = 1
$ws.Range("O7").Value = 16
$ws.Range("F8").Value = 45688.43347222222
$ws.Range("O8").Value = 12
$ws.Range("F9").Value = 45688.43751157408
$ws.Range("O9").Value = 18
$ws.Range("F10").Value = 45688.43893518519
$ws.Range("O10").Value = 20
$ws.Range("F11").Value = 45688.43164351852
$ws.Range("O11").Value = 10
$ws.Range("F12").Value = 45688.43100694445
$ws.Range("O12").Value = 9
$ws.Range("F13").Value = 45688.44256944444
$ws.Range("O13").Value = 26
$ws.Range("F14").Value = 45688.43516203704
$ws.Range("O14").Value = 15
$ws.Range("F15").Value = 45688.44451388889
$ws.Range("O15").Value = 28

$ws = $wb.Worksheets.Item("NextBus3")
$ws.Range("O2").Value = 33
$ws.Range("F3").Value = 45688.44222222222
$ws.Range("O3").Value = 25
$ws.Range("F4").Value = 45688.43768518518
$ws.Range("L4").Value = "BD"
$ws.Range("O4").Value = 19
$ws.Range("F5").Value = 45688.45008101852
$ws.Range("O5").Value = 36
$ws.Range("F6").Value = 45688.4377662037
$ws.Range("O6").Value = 19
$ws.Range("O7").Value = 26
$ws.Range("F8").Value = 45688.43837962963
$ws.Range("O8").Value = 20
$ws.Range("F9").Value = 45688.44399305555
$ws.Range("O9").Value = 28
$ws.Range("F10").Value = 45688.45055555556
$ws.Range("O10").Value = 37
$ws.Range("O11").Value = 18
$ws.Range("F12").Value = 45688.43612268518
$ws.Range("O12").Value = 16
$ws.Range("F13").Value = 45688.44969907407
$ws.Range("O13").Value = 36
$ws.Range("O14").Value = 21
$ws.Range("O15").Value = 41
